# Refresh market-data columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# on the Aegis_Profits leve-profit sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 640.05
$arr[0,1] = 139.33333
$arr[0,2] = 1391.125
$arr[0,3] = 139.33333
$arr[0,4] = 1391.125
$arr[0,5] = 89.66667000000001
$arr[0,6] = -1849.125
$ws.Range("H33:N33").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7828.8335
$arr[0,1] = 6753.25
$arr[0,2] = 9980
$arr[0,3] = 20259.75
$arr[0,4] = 29940
$arr[0,5] = -19385.75
$arr[0,6] = -31688
$ws.Range("H69:N69").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 7828.8335
$arr[0,1] = 6753.25
$arr[0,2] = 9980
$arr[0,3] = 60779.25
$arr[0,4] = 89820
$arr[0,5] = -56411.25
$arr[0,6] = -98556
$ws.Range("H72:N72").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1179.4517
$arr[0,1] = 798.4583
$arr[0,2] = 2485.7144
$arr[0,3] = 798.4583
$arr[0,4] = 2485.7144
$arr[0,5] = 449.5417
$arr[0,6] = -4981.7144
$ws.Range("H92:N92").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2223.9285
$arr[0,1] = 1757.2222
$arr[0,2] = 3064
$arr[0,3] = 1757.2222
$arr[0,4] = 3064
$arr[0,5] = -1216.2222
$arr[0,6] = -4146
$ws.Range("H100:N100").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 68348
$arr[0,1] = 144271.42
$arr[0,2] = 1915
$arr[0,3] = 144271.42
$arr[0,4] = 1915
$arr[0,5] = -141017.42
$arr[0,6] = -8423
$ws.Range("H113:N113").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1118.875
$arr[0,1] = 1197.9166
$arr[0,2] = 881.75
$arr[0,3] = 10781.2494
$arr[0,4] = 7935.75
$arr[0,5] = -8246.2494
$arr[0,6] = -13005.75
$ws.Range("H135:N135").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6411.875
$arr[0,1] = 10500
$arr[0,2] = 5827.857
$arr[0,3] = 31500
$arr[0,4] = 17483.571
$arr[0,5] = -26320
$arr[0,6] = -27843.571
$ws.Range("H141:N141").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5353.727
$arr[0,1] = 1722.75
$arr[0,2] = 7428.5713
$arr[0,3] = 1722.75
$arr[0,4] = 7428.5713
$arr[0,5] = -1607.75
$arr[0,6] = -7658.5713
$ws.Range("H3:N3").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 20140.197
$arr[0,1] = 3583.3972
$arr[0,2] = 87287.22
$arr[0,3] = 3583.3972
$arr[0,4] = 87287.22
$arr[0,5] = -3296.3972
$arr[0,6] = -87861.22
$ws.Range("H32:N32").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1842.6279
$arr[0,1] = 1106.091
$arr[0,2] = 2614.238
$arr[0,3] = 1106.091
$arr[0,4] = 2614.238
$arr[0,5] = -894.0909999999999
$arr[0,6] = -3038.238
$ws.Range("H61:N61").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 79522.305
$arr[0,1] = 126583.625
$arr[0,2] = 4224.2
$arr[0,3] = 126583.625
$arr[0,4] = 4224.2
$arr[0,5] = -126087.625
$arr[0,6] = -5216.2
$ws.Range("H97:N97").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1842.6279
$arr[0,1] = 1106.091
$arr[0,2] = 2614.238
$arr[0,3] = 3318.273
$arr[0,4] = 7842.714
$arr[0,5] = -768.2729999999997
$arr[0,6] = -12942.714
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 21627
$arr[0,1] = 30115.766
$arr[0,2] = 1011.4286
$arr[0,3] = 30115.766
$arr[0,4] = 1011.4286
$arr[0,5] = -29664.766
$arr[0,6] = -1913.4286
$ws.Range("H94:N94").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 75726.37
$arr[0,1] = 84763
$arr[0,2] = 68497.07
$arr[0,3] = 84763
$arr[0,4] = 68497.07
$arr[0,5] = -83016
$arr[0,6] = -71991.07
$ws.Range("H105:N105").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1638.8649
$arr[0,1] = 1507.5758
$arr[0,2] = 2722
$arr[0,3] = 4522.7274
$arr[0,4] = 8166
$arr[0,5] = -1987.7274
$arr[0,6] = -13236
$ws.Range("H134:N134").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3200.2
$arr[0,1] = 3302
$arr[0,2] = 3174.75
$arr[0,3] = 3302
$arr[0,4] = 3174.75
$arr[0,5] = -2678
$arr[0,6] = -4422.75
$ws.Range("H62:N62").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3200.2
$arr[0,1] = 3302
$arr[0,2] = 3174.75
$arr[0,3] = 16510
$arr[0,4] = 15873.75
$arr[0,5] = -13390
$arr[0,6] = -22113.75
$ws.Range("H65:N65").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1008.0244
$arr[0,1] = 750.2059
$arr[0,2] = 2260.2856
$arr[0,3] = 2250.6177
$arr[0,4] = 6780.8568
$arr[0,5] = 284.3822999999998
$arr[0,6] = -11850.8568
$ws.Range("H134:N134").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 16610.738
$arr[0,1] = 943.95
$arr[0,2] = 23573.756
$arr[0,3] = 2831.85
$arr[0,4] = 70721.26800000001
$arr[0,5] = -2020.85
$arr[0,6] = -72343.26800000001
$ws.Range("H68:N68").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 16610.738
$arr[0,1] = 943.95
$arr[0,2] = 23573.756
$arr[0,3] = 8495.550000000001
$arr[0,4] = 212163.804
$arr[0,5] = -4439.550000000001
$arr[0,6] = -220275.804
$ws.Range("H71:N71").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3130.8333
$arr[0,1] = 4201.6665
$arr[0,2] = 2060
$arr[0,3] = 4201.6665
$arr[0,4] = 2060
$arr[0,5] = -3203.6665
$arr[0,6] = -4056
$ws.Range("H80:N80").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3130.8333
$arr[0,1] = 4201.6665
$arr[0,2] = 2060
$arr[0,3] = 21008.3325
$arr[0,4] = 10300
$arr[0,5] = -16016.3325
$arr[0,6] = -20284
$ws.Range("H83:N83").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2233.3333
$arr[0,1] = 2250
$arr[0,2] = 2200
$arr[0,3] = 2250
$arr[0,4] = 2200
$arr[0,5] = -1709
$arr[0,6] = -3282
$ws.Range("H100:N100").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 222856.67
$arr[0,1] = 250532.75
$arr[0,2] = 200715.8
$arr[0,3] = 501065.5
$arr[0,4] = 401431.6
$arr[0,5] = -500004.5
$arr[0,6] = -403553.6
$ws.Range("H81:N81").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 222856.67
$arr[0,1] = 250532.75
$arr[0,2] = 200715.8
$arr[0,3] = 2505327.5
$arr[0,4] = 2007158
$arr[0,5] = -2500023.5
$arr[0,6] = -2017766
$ws.Range("H84:N84").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 111639.22
$arr[0,1] = 167068.83
$arr[0,2] = 780
$arr[0,3] = 334137.66
$arr[0,4] = 1560
$arr[0,5] = -333596.66
$arr[0,6] = -2642
$ws.Range("H100:N100").Value = $arr

# CUL rows where a cell was inserted or removed mid-row (column layout shifts)
$ws = $wb.Worksheets.Item("CUL")

# Row 6: N6 is removed; M6 becomes the final (combined) value
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 112.111115
$arr[0,1] = 112.111115
$arr[0,2] = 0
$arr[0,3] = 336.333345
$arr[0,4] = 0
$ws.Range("H6:L6").Value = $arr
$ws.Range("M6").Value = -223.333345
$ws.Range("N6").ClearContents()

# Row 32: M32 is removed; N32 becomes the final value
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1912.8572
$arr[0,1] = 0
$arr[0,2] = 1912.8572
$arr[0,3] = 0
$ws.Range("H32:K32").Value = $arr
$ws.Range("L32").Value = 5738.571599999999
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -6304.571599999999

# Row 34: M34 is newly added (row previously skipped straight from L to N)
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 1975
$arr[0,1] = 1000
$arr[0,2] = 2462.5
$arr[0,3] = 3000
$ws.Range("H34:K34").Value = $arr
$ws.Range("L34").Value = 7387.5
$ws.Range("M34").Value = -2916
$ws.Range("N34").Value = -7555.5

# Row 131: M131 is newly added (row previously skipped straight from L to N)
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 847.11
$arr[0,1] = 599
$arr[0,2] = 849.61615
$arr[0,3] = 1797
$ws.Range("H131:K131").Value = $arr
$ws.Range("L131").Value = 2548.84845
$ws.Range("M131").Value = 3243
$ws.Range("N131").Value = -12628.84845

Write-Output "Applied Aegis_Profits market-data refresh"